# Generate Report for Handoff
# The file "9cb7d7da-ca99-407d-8a7f-a2d0a2698a16.md" has moved from "In
# Translation" to "Ready for handoff" status, with a machine-translation
# ("mt") priority and refreshed handoff timestamps, across all three
# worksheets (Overview, zh-cn, de-de).

$wb = $excel.ActiveWorkbook

# --- Overview sheet: row for 9cb7d7da-ca99-407d-8a7f-a2d0a2698a16.md is row 3 ---
$wsOverview = $wb.Worksheets.Item("Overview")
$wsOverview.Range("E3").Value = "Ready for handoff"
$wsOverview.Range("F3").Value = "Ready for handoff"
$wsOverview.Range("G3").Value = "2016-08-30 10:14:03"

# --- zh-cn sheet: row for 9cb7d7da-ca99-407d-8a7f-a2d0a2698a16.md is row 3 ---
$wsZhCn = $wb.Worksheets.Item("zh-cn")
$wsZhCn.Range("C3").Value = "Ready for handoff"
$wsZhCn.Range("E3").Value = "mt"
$wsZhCn.Range("H3").Value = "2016-08-30 10:13:56"

# --- de-de sheet: row for 9cb7d7da-ca99-407d-8a7f-a2d0a2698a16.md is row 3 ---
$wsDeDe = $wb.Worksheets.Item("de-de")
$wsDeDe.Range("C3").Value = "Ready for handoff"
$wsDeDe.Range("E3").Value = "mt"
$wsDeDe.Range("H3").Value = "2016-08-30 10:14:03"

# The longer "Ready for handoff" status text widens its column (Status on
# zh-cn/de-de; zh-cn/de-de on Overview) the same way Excel auto-fits a
# column after a cell's content changes.
$wsOverview.Columns.Item(5).ColumnWidth = 16.37
$wsOverview.Columns.Item(6).ColumnWidth = 16.37
$wsZhCn.Columns.Item(3).ColumnWidth = 16.37
$wsDeDe.Columns.Item(3).ColumnWidth = 16.37
